# "cw done final commit"
# Reset the player standings: TOTAL RUNS (col B) and WICKETS (col C) go back to 0
# for every player, then the real (small) tallies that are actually known get
# filled back in for the rows that have them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero out the whole data range (rows 2-89, columns B:C) first.
$ws.Range("B2:C89").Value = 0

# Now write back the non-zero totals/wickets that survived the reset.
$nonZero = @{
    "B57" = 55;  "B58" = 30;  "B59" = 7;   "B60" = 18;  "B61" = 14;  "B62" = 36
    "B63" = 8;   "B64" = 3;   "C64" = 1
    "C65" = 3
    "C66" = 2
    "C67" = 4
    "B79" = 10;  "B80" = 18;  "B81" = 27;  "B82" = 21;  "B83" = 5;   "B84" = 55
    "B85" = 3;   "C85" = 2
    "B86" = 8;   "C86" = 1
    "B87" = 11;  "C87" = 1
    "B88" = 9;   "C88" = 1
    "B89" = 1;   "C89" = 1
}

foreach ($addr in $nonZero.Keys) {
    $ws.Range($addr).Value = $nonZero[$addr]
}

# Match the author's final selection (B2:C89, anchored at B2).
$ws.Range("B2:C89").Select()
